$wb = $excel.ActiveWorkbook

# 1. Add a new blank sheet ("Sheet4") at the end of the workbook.
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws4.Name = "Sheet4"

# 2. New project mock-up data on Sheet2 (Amount column + two more data rows).
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("H19").Value = "Amount"
$ws2.Range("H20").Value = 24000
$ws2.Range("H21").Value = 35000
$ws2.Range("I21").Value = 59000

# 3. Populate Sheet3 with the loan-request "Stages" / "Document Checklist" / "Folder Structure" mock-up.
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("H9").Value = "Stages"
$ws3.Range("N9").Value = "Document Checklist"
$ws3.Range("Q9").Value = "Folder Structure"
$ws3.Range("H11").Value = "Active Pre Lead"
$ws3.Range("I11").Value = "Borrower makes request."
$ws3.Range("N11").Value = "KYC Documents"
$ws3.Range("Q11").Value = "Confidential"
$ws3.Range("H12").Value = "Pre Lead Client Drop"
$ws3.Range("I12").Value = "Borrower or Systme Admin deletes loan request."
$ws3.Range("N12").Value = "Certificate of Incorporation"
$ws3.Range("Q12").Value = "Public"
$ws3.Range("H13").Value = "Active Stage 1 Lead"
$ws3.Range("I13").Value = "Background verification done Mandate letter generated"
$ws3.Range("N13").Value = "MOA and AOA"
$ws3.Range("Q13").Value = "Subfolders under Public"
$ws3.Range("H14").Value = "Stage 1 Reject"
$ws3.Range("I14").Value = "Sales rep. Rejects loan request at time of background verification."
$ws3.Range("N14").Value = "Pan Card of Company"
$ws3.Range("Q14").Value = "KYC"
$ws3.Range("H15").Value = "Stage 1 Client Drop"
$ws3.Range("I15").Value = "Borrower or Systme Admin deletes loan request."
$ws3.Range("N15").Value = "List of Directors"
$ws3.Range("Q15").Value = "Financials"
$ws3.Range("H16").Value = "Active Stage 2 Lead"
$ws3.Range("I16").Value = "Borrower upload document after successful background check"
$ws3.Range("N16").Value = "Pan Cards of all directors"
$ws3.Range("Q16").Value = "Tax Returns"
$ws3.Range("H17").Value = "Stage 2 Reject"
$ws3.Range("I17").Value = "Credit checker reject loan req."
$ws3.Range("N17").Value = "Address Proof of all directors"
$ws3.Range("Q17").Value = "Bank Statements"
$ws3.Range("H18").Value = "Stage 2 Client Drop"
$ws3.Range("I18").Value = "Borrower or Systme Admin deletes loan request."
$ws3.Range("N18").Value = "Financial Statements"
$ws3.Range("Q18").Value = "PO and Invoice"
$ws3.Range("H19").Value = "Active Listing"
$ws3.Range("I19").Value = "Credit checker approves loan request and loan request is forwareded to Lenders"
$ws3.Range("N19").Value = "Audited financials for past three years"
$ws3.Range("Q19").Value = "Other"
$ws3.Range("H20").Value = "Post Listing Client Drop"
$ws3.Range("I20").Value = "Borrower or Systme Admin deletes loan request."
$ws3.Range("N20").Value = "Audit Report and Tax Audit Report for last 3 years"
$ws3.Range("Q20").Value = "Lender 1 Additional Info"
$ws3.Range("H21").Value = "Lender Reject"
$ws3.Range("I21").Value = "Lender rejects"
$ws3.Range("N21").Value = "Provisional financials for current year"
$ws3.Range("Q21").Value = "Lender 2 Additional Info"
$ws3.Range("H22").Value = "Active Sanction"
$ws3.Range("I22").Value = "Lender senction loan request"
$ws3.Range("N22").Value = "Projected financials for loan term"
$ws3.Range("Q22").Value = "Lender 3 Additional Info"
$ws3.Range("H23").Value = "Sanctioned Client Drop"
$ws3.Range("I23").Value = "Borrower or Systme Admin deletes loan request."
$ws3.Range("N23").Value = "Past Investments"
$ws3.Range("H24").Value = "Sanctioned Lender Drop"
$ws3.Range("I24").Value = "??"
$ws3.Range("N24").Value = "Shareholding Pattern"
$ws3.Range("H25").Value = "Disbursed"
$ws3.Range("I25").Value = "??"
$ws3.Range("N25").Value = "Year-wise breakup of equity investment received (amount and investor)"
$ws3.Range("H26").Value = "Active Repeat Listing"
$ws3.Range("I26").Value = "??"
$ws3.Range("N26").Value = "Loans and Banking"
$ws3.Range("N27").Value = "List of all bank accounts maintained by the company"
$ws3.Range("N28").Value = "Bank Statements for ALL bank accounts for the last 12 months"
$ws3.Range("N29").Value = "Term sheets for all loans and credit facilities currently outstanding"
$ws3.Range("N30").Value = "Business"
$ws3.Range("N31").Value = "Company profile and brochure, if any"
$ws3.Range("N32").Value = "Latest investor presentation deck"
$ws3.Range("N33").Value = "Organisation Structure"
$ws3.Range("N34").Value = "Statutory Returns"
$ws3.Range("N35").Value = "Income Tax Returns for the last 2 years"
$ws3.Range("N36").Value = "VAT and/or service tax returns for the last 12 months"
$ws3.Range("N37").Value = "For Unsecured Business Loans or Invoice Discounting"
$ws3.Range("N38").Value = "List of Top Customers (Top 10 customers or top 50% sales)"
$ws3.Range("N39").Value = "Debtor Ageing Schedule"
$ws3.Range("N40").Value = "Customer contracts"
$ws3.Range("N41").Value = "Customer ledgers for top 10 customers"
$ws3.Range("N42").Value = "For PO based loan"
$ws3.Range("N43").Value = "Copy of the contract"
$ws3.Range("N44").Value = "Tender on which basis the contract was awarded"
$ws3.Range("N45").Value = "Milestones and expected payments/receipts for the project"
$ws3.Range("N46").Value = "Case studies - similar projects executed in the past"
$ws3.Range("N47").Value = "For term loan and leasing"
$ws3.Range("N48").Value = "Breakup of project cost"
$ws3.Range("N49").Value = "List of assets to be purchased"
$ws3.Range("N50").Value = "Technical specifications of proposed project"
$ws3.Range("N51").Value = "List of suppliers"
$ws3.Range("N52").Value = "Cash flow projections for the project"

# 4. Restore selections / active sheet so the saved view matches the authored state.
$ws2.Activate()
$ws2.Range("F11").Select()
$ws3.Range("Q9").Select()
$ws4.Range("E5").Select()
$ws2.Activate()
